$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (sheet name used in workbook.xml)
$ws.Name = "876987"

# Delete the old header row and the first two data rows (省份/销量, 广东, 广西),
# shifting the remaining rows up by three.
$ws.Range("A1:B3").EntireRow.Delete()

# Clear out column B entirely for the remaining rows -- the cells stay present
# (so the sheet's used range still spans column B) but hold no value, matching
# the "t=str" empty-value cells in the target sheet.
$ws.Range("B1:B12").ClearContents()
$ws.Range("B1:B12").Style = "Normal"
